# Updated cryptos list with latest price/volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.521.78'
$ws.Range('E2').Value = '  +1.94%  '
$ws.Range('D3').Value = '2.547.18'
$ws.Range('E3').Value = '  +4.69%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''571.98'
$ws.Range('E5').Value = '  +2.77%  '
$ws.Range('D6').Value = '''151.10'
$ws.Range('E6').Value = '  +8.88%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '''0.589'
$ws.Range('E8').Value = '  +0.94%  '
$ws.Range('D9').Value = '2.545.45'
$ws.Range('E9').Value = '  +4.67%  '
$ws.Range('E10').Value = '  +2.50%  '
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('D14').Value = '''28.36'
$ws.Range('E14').Value = '  +8.30%  '
$ws.Range('D15').Value = '3.003.31'
$ws.Range('E15').Value = '  +4.74%  '
$ws.Range('D16').Value = '63.435.76'
$ws.Range('E16').Value = '  +2.08%  '
$ws.Range('E17').Value = '  +1.44%  '
$ws.Range('D18').Value = '2.557.90'
$ws.Range('E18').Value = '  +5.33%  '
$ws.Range('E19').Value = '  +4.22%  '
$ws.Range('D20').Value = '''340.96'
$ws.Range('E20').Value = '  -1.52%  '
$ws.Range('E21').Value = '  +4.25%  '
$ws.Range('E22').Value = '  +1.23%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = '''66.21'
$ws.Range('E24').Value = '  +1.54%  '
$ws.Range('E25').Value = '  -1.05%  '
$ws.Range('D26').Value = '''1.61'
$ws.Range('E26').Value = '  +5.38%  '
$ws.Range('D27').Value = '''1.50'
$ws.Range('E27').Value = '  +12.67%  '
$ws.Range('D28').Value = '''8.53'
$ws.Range('E28').Value = '  +3.90%  '
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('D30').Value = '''7.13'
$ws.Range('E30').Value = '  +12.01%  '
$ws.Range('D31').Value = '0.0₃0837'
$ws.Range('E31').Value = '  +6.86%  '
$ws.Range('E32').Value = '  +3.63%  '
$ws.Range('D33').Value = '''177.98'
$ws.Range('E33').Value = '  +3.59%  '
$ws.Range('E34').Value = '  +9.35%  '
$ws.Range('D35').Value = '''420.96'
$ws.Range('E35').Value = '  +14.96%  '
$ws.Range('D36').Value = '''0.407'
$ws.Range('E36').Value = '  +2.76%  '
$ws.Range('D37').Value = '''19.23'
$ws.Range('E37').Value = '  +3.40%  '
$ws.Range('D38').Value = '''4.46'
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('E40').Value = '  +4.92%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').Value = '''40.03'
$ws.Range('E42').Value = '  +2.13%  '
$ws.Range('D43').Value = '''154.84'
$ws.Range('E43').Value = '  +5.95%  '
$ws.Range('E44').Value = '  +4.12%  '
$ws.Range('D45').Value = '''21.18'
$ws.Range('E45').Value = '  +3.17%  '
$ws.Range('D46').Value = '''0.613'
$ws.Range('E46').Value = '  +4.26%  '
$ws.Range('E47').Value = '  +2.89%  '
$ws.Range('E48').Value = '  +9.40%  '
$ws.Range('D49').Value = '''0.0970'
$ws.Range('E49').Value = '  +1.74%  '
$ws.Range('D50').Value = '''18.63'
$ws.Range('E50').Value = '  +4.25%  '
$ws.Range('E51').Value = '  +7.51%  '
